$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A26").Value = "Camping"
$ws.Range("A13").Value = "Albergue Costa del Caribe"
$ws.Range("A14").Value = "Albergue EcoAlbergue Sierra Verde"

$ws.Activate()
$ws.Range("E12").Select()
$excel.ActiveWindow.Zoom = 57
